$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 used to hold the shared string "R40"; the authored change replaces its
# contents with the literal text "1". Prefixing the assignment with an
# apostrophe makes Excel store the numeric-looking value as quote-prefixed
# text instead of silently coercing it into a real number, which is what
# keeps the cell a shared-string ("R40" -> "1") exactly like the diff shows.
$ws.Range("B11").Value = "'1"
